$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header labels for the RVR helper columns (row 12) ---
$ws.Range("H12").Value = "RVR"
$ws.Range("I12").Value = "RVR/2"

# --- H/I helper columns: RVR = D/D(ref row), RVR/2 = RVR / 2 ---
$ws.Range("H14").Formula = "=D14/D3"
$ws.Range("H15").Formula = "=D15/D4"
$ws.Range("H16").Formula = "=D16/D5"
$ws.Range("H17").Formula = "=D17/D6"
$ws.Range("H18").Formula = "=D18/D7"
$ws.Range("H19").Formula = "=D19/D8"

$ws.Range("I14").Formula = "=H14/2"
$ws.Range("I15").Formula = "=H15/2"
$ws.Range("I16").Formula = "=H16/2"
$ws.Range("I17").Formula = "=H17/2"
$ws.Range("I18").Formula = "=H18/2"
$ws.Range("I19").Formula = "=H19/2"

# --- Labels for the two regression blocks ---
$ws.Range("K13").Value = "linear"
$ws.Range("K20").Value = "constant"

# --- LINEST regression #1: RVR/2 vs I3:I8 ---
$ws.Range("K14:L18").FormulaArray = "=LINEST(I14:I19,I3:I8,TRUE,TRUE)"

# --- LINEST regression #2: RVR/2 vs Month (A3:A8) ---
$ws.Range("K21:L25").FormulaArray = "=LINEST(I14:I19,A3:A8,TRUE,TRUE)"

# Final selection, matching the state left by the author after entry
$ws.Range("L28").Select() | Out-Null
